$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-converts decimal-looking text (e.g. "1.000", "0.7212") into a
# real number when assigned via .Value, which both changes the cell's type
# (numeric instead of text) and mangles the literal text (trailing zeros /
# floating point noise). Force those through as text: mark the cell Text
# ("@") before writing the literal, then clear the number-format stamp so
# the on-disk cell keeps its original (unstyled) look while the stored
# value stays a verbatim string.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '30.287.32'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.931.02'
$ws.Range("E3").Value = '  -0.11%  '
Set-TextValue "D4" '1.001'
$ws.Range("E4").Value = '  +0.09%  '
Set-TextValue "D5" '249.72'
$ws.Range("E5").Value = '  +0.29%  '
Set-TextValue "D6" '0.7212'
$ws.Range("E6").Value = '  +1.29%  '
Set-TextValue "D7" '1.000'
$ws.Range("E7").Value = '  +0.01%  '
Set-TextValue "D8" '0.3215'
$ws.Range("E8").Value = '  -2.00%  '
Set-TextValue "D9" '27.61'
$ws.Range("E9").Value = '  -0.10%  '
Set-TextValue "D10" '0.07080'
$ws.Range("E10").Value = '  +3.30%  '
Set-TextValue "D11" '0.7912'
$ws.Range("E11").Value = '  -1.97%  '
Set-TextValue "D12" '0.08042'
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").Value = '1.928.73'
$ws.Range("E13").Value = '  -0.19%  '
Set-TextValue "D14" '5.380'
$ws.Range("E14").Value = '  -1.22%  '
Set-TextValue "D15" '94.87'
$ws.Range("E15").Value = '  -0.22%  '
Set-TextValue "D16" '14.61'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '30.303.68'
$ws.Range("E17").Value = '  -0.07%  '
Set-TextValue "D18" '258.29'
$ws.Range("E18").Value = '  -0.63%  '
Set-TextValue "D19" '0.000008041'
$ws.Range("E19").Value = '  -0.06%  '
Set-TextValue "D20" '5.735'
$ws.Range("E20").Value = '  -1.95%  '
$ws.Range("D21").Value = '2.185.40'
$ws.Range("E21").Value = '  -0.03%  '
Set-TextValue "D22" '1.000'
$ws.Range("E22").Value = '  +0.03%  '
Set-TextValue "D23" '1.002'
$ws.Range("E23").Value = '  +0.22%  '
Set-TextValue "D24" '6.829'
$ws.Range("E24").Value = '  -0.78%  '
Set-TextValue "D25" '9.529'
$ws.Range("E25").Value = '  -2.05%  '
Set-TextValue "D26" '165.37'
$ws.Range("E26").Value = '  +3.54%  '
Set-TextValue "D27" '19.26'
$ws.Range("E27").Value = '  +0.56%  '
Set-TextValue "D28" '2.288'
$ws.Range("E28").Value = '  -3.55%  '
Set-TextValue "D29" '0.1285'
$ws.Range("E29").Value = '  -3.55%  '
Set-TextValue "D30" '1.367'
$ws.Range("E30").Value = '  +1.29%  '
Set-TextValue "D31" '1.534'
$ws.Range("E31").Value = '  -1.75%  '
Set-TextValue "D32" '4.403'
$ws.Range("E32").Value = '  -0.34%  '
Set-TextValue "D33" '4.159'
$ws.Range("E33").Value = '  -1.50%  '
Set-TextValue "D34" '0.05193'
$ws.Range("E34").Value = '  +1.92%  '
Set-TextValue "D35" '1.259'
$ws.Range("E35").Value = '  +2.86%  '
Set-TextValue "D36" '0.7448'
$ws.Range("E36").Value = '  +0.26%  '
Set-TextValue "D37" '2.770'
$ws.Range("E37").Value = '  +0.09%  '
Set-TextValue "D38" '0.01956'
$ws.Range("E38").Value = '  -1.30%  '
Set-TextValue "D39" '2.810'
$ws.Range("E39").Value = '  -0.34%  '
Set-TextValue "D40" '77.61'
$ws.Range("E40").Value = '  -2.39%  '
Set-TextValue "D41" '6.383'
$ws.Range("E41").Value = '  -3.35%  '
Set-TextValue "D42" '0.4501'
$ws.Range("E42").Value = '  +0.57%  '
Set-TextValue "D43" '2.001'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D44" '1.000'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D45" '0.8393'
$ws.Range("E45").Value = '  +0.37%  '
Set-TextValue "D46" '100.78'
$ws.Range("E46").Value = '  -1.25%  '
Set-TextValue "D47" '9.775'
$ws.Range("E47").Value = '  -0.04%  '
Set-TextValue "D48" '7.443'
$ws.Range("E48").Value = '  +1.54%  '
Set-TextValue "D49" '36.45'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D50" '0.06102'
$ws.Range("E50").Value = '  +2.58%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue "D51" '0.4178'
$ws.Range("E51").Value = '  +2.07%  '

Write-Output "Updated cryptos list"
